$d = $word.ActiveDocument

# --- Block A: JSON sample values list (Courier New, sz16) ---
# "frameSize": 1016,            -> "eepromFrameCount": 256,
# "frameSizeWithCRC": 1024,     -> "eepromFramePayloadSize": 1016,
# "framesCount": 256,           -> "eepromFrameSize": 1024,
$d.Content.Find.Execute('"frameSize": 1016,', $true, $false, $false, $false, $false, `
    $true, 1, $false, '"eepromFrameCount": 256,', 2) | Out-Null

$d.Content.Find.Execute('"frameSizeWithCRC": 1024,', $true, $false, $false, $false, $false, `
    $true, 1, $false, '"eepromFramePayloadSize": 1016,', 2) | Out-Null

$d.Content.Find.Execute('"framesCount": 256,', $true, $false, $false, $false, $false, `
    $true, 1, $false, '"eepromFrameSize": 1024,', 2) | Out-Null

# --- Block B: descriptive bullet list (Times New Roman, sz28) ---
# - "frameSize"-size of data frame, in bytes;
# - "frameSizeWithCRC"-size of data frame with CRC, in bytes;
# - "framesCount"-number of data frames in the file;
$d.Content.Find.Execute([char]0x2013 + ' ' + [char]0x201C + 'frameSize' + [char]0x201D + [char]0x2013 + 'size of data frame, in bytes;', `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    [char]0x2013 + ' ' + [char]0x201C + 'eepromFramePayloadSize' + [char]0x201D + [char]0x2013 + 'size of data frame, in bytes;', 2) | Out-Null

$d.Content.Find.Execute([char]0x2013 + ' ' + [char]0x201C + 'frameSizeWithCRC' + [char]0x201D + [char]0x2013 + 'size of data frame with CRC, in bytes;', `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    [char]0x2013 + ' ' + [char]0x201C + 'eepromFrameSize' + [char]0x201D + [char]0x2013 + 'size of data frame with CRC, in bytes;', 2) | Out-Null

$d.Content.Find.Execute([char]0x2013 + ' ' + [char]0x201C + 'framesCount' + [char]0x201D + [char]0x2013 + 'number of data frames in the file;', `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    [char]0x2013 + ' ' + [char]0x201C + 'eepromFrameCount' + [char]0x201D + [char]0x2013 + 'number of data frames in the file;', 2) | Out-Null

Write-Host "Done."
